$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark: remove it from its old spot (between
#    "You have" + " " and "{{ nice_number(...) }}") in the "Print, Copy
#    and Deliver" section.
# ---------------------------------------------------------------------
$find = $d.Content
$found = $find.Find.Execute("You have", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the 'You have' paragraph to relocate the _GoBack bookmark."
}
$oldPara = $find.Paragraphs(1)
$oldParaRange = $oldPara.Range

$oldParaXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="63E59D45" w14:textId="43AADAB8" w:rsidR="00976B3B" w:rsidRDefault="00976B3B" w:rsidP="00CA3FEC"><w:r><w:t>You have</w:t></w:r><w:r w:rsidR="00C41FCE"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="0091371C"><w:t xml:space="preserve">{{ </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="0091371C"><w:t>nice_number</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="0091371C"><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="0091371C"><w:t>num_downloads</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="0091371C"><w:t>) }}</w:t></w:r><w:r w:rsidR="008C29F4"><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>forms. Make sure to print, copy, and organize 3 sets of each form.</w:t></w:r></w:p>
'@

$oldParaRange.InsertXML($oldParaXml) | Out-Null

# ---------------------------------------------------------------------
# 2) At the end of the document (after the "Arrive by 9:00 AM..."
#    paragraph) add a blank paragraph, then an italicized, small (10pt)
#    help paragraph that reports the guided-interview's update date and
#    the forms' generation date. Put the "_GoBack" bookmark at the very
#    end of that new paragraph.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$lastParaRange = $lastPara.Range
$insertionPoint = $d.Range($lastParaRange.End - 1, $lastParaRange.End - 1)

$newParasXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:i/><w:sz w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:sz w:val="20"/></w:rPr><w:t>Guided i</w:t></w:r><w:r><w:rPr><w:i/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">nterview </w:t></w:r><w:r><w:rPr><w:i/><w:sz w:val="20"/></w:rPr><w:t>update date</w:t></w:r><w:r><w:rPr><w:i/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve"> is {{ </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/><w:sz w:val="20"/></w:rPr><w:t>format_date</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/><w:sz w:val="20"/></w:rPr><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/><w:sz w:val="20"/></w:rPr><w:t>all_variables</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/><w:sz w:val="20"/></w:rPr><w:t>(special='metadata')['</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/><w:sz w:val="20"/></w:rPr><w:t>revision_date</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/><w:sz w:val="20"/></w:rPr><w:t>']) }}</w:t></w:r><w:r><w:rPr><w:i/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve"> and these forms were generated {{ today() }}.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

$insertionPoint.InsertXML($newParasXml) | Out-Null
